$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new shared string content "get file quickly" in C3
$ws.Range("C3").Value = "get file quickly"

# Update selection to C6 (does not need to contain data)
$ws.Range("C6").Select()
